# Update the "Generate Report for Handback" timestamps in the handback-status workbook.
$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for a80a1b69-...
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 22:57:09"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for a80a1b69-...
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-05 22:56:59"
$wsZhCn.Range("K4").Value = "2016-09-05 22:57:39"

# de-de sheet: Correspond Handback DateTime for a80a1b69-...
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-05 22:57:46"
